$wb = $excel.ActiveWorkbook

$sheetNames = @("V1", "V2", "V3", "V4")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # New "Request/sec" column next to the first (Summary) table's Throughput column (H)
    $ws.Range("L2").Value = "Request/sec"
    $ws.Range("L3").Formula = "=H3*60"
    $ws.Range("L4").Formula = "=H4*60"

    # New "Request/sec" column next to the second (Aggregate) table's Throughput column (K)
    $ws.Range("N9").Value = "Request/sec"
    $ws.Range("N10").Formula = "=K10*60"
    $ws.Range("N11").Formula = "=K11*60"

    # Mirror the author's manual selection of the new column while filling it in
    $ws.Range("L2:L4").Select() | Out-Null
}

# Leave V1 as the active sheet/tab, matching the final saved view state
$wb.Worksheets.Item("V1").Activate() | Out-Null
